$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Copy formatting from row 5 down into new rows 6 and 7 first
$ws.Range("A5:F5").Copy()
$ws.Range("A6:F6").PasteSpecial($xlPasteFormats)
$ws.Range("A5:F5").Copy()
$ws.Range("A7:F7").PasteSpecial($xlPasteFormats)

# Row 6: 2025-10-03, 四方坪站
$ws.Cells.Item(6, 1).Value = 45933
$ws.Cells.Item(6, 2).Value = "四方坪站"
$ws.Cells.Item(6, 3).Value = 9382.7900000000009
$ws.Cells.Item(6, 4).Value = 7814.8
$ws.Cells.Item(6, 5).Value = 3200.49
$ws.Cells.Item(6, 6).Value = 370

# Row 7: 2025-10-03, 高岭站
$ws.Cells.Item(7, 1).Value = 45933
$ws.Cells.Item(7, 2).Value = "高岭站"
$ws.Cells.Item(7, 3).Value = 3397.12
$ws.Cells.Item(7, 4).Value = 2536.4
$ws.Cells.Item(7, 5).Value = 850.53
$ws.Cells.Item(7, 6).Value = 116

# Update selection to match target
$ws.Range("J11").Select()
